$d = $word.ActiveDocument

# Locate the end of the "Dewey, J. (1904). Ethics." sentence.
$find = $d.Content
$found = $find.Find.Execute("Dewey, J. (1904). Ethics.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertionPoint = $find.Duplicate
    $insertionPoint.Collapse(0)

    # Track the new text as discrete insertions (and then accept them
    # individually) so each piece of text lands in its own <w:r> run,
    # matching how the citation suffix was actually typed in, instead of
    # being silently coalesced into a single run with the preceding text.
    $wasTracking = $d.TrackRevisions
    $d.TrackRevisions = $true

    $r1 = $d.Range($insertionPoint.Start, $insertionPoint.Start)
    $r1.InsertAfter("(MW3: ")

    $r2 = $d.Range($r1.End, $r1.End)
    $r2.InsertAfter("40-")

    $r3 = $d.Range($r2.End, $r2.End)
    $r3.InsertAfter("58). ")

    $d.TrackRevisions = $wasTracking

    for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
        $d.Revisions.Item($i).Accept()
    }
}
